$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = [double]"25.14000000000049"
$ws.Range("G2").Value = [double]"1.110223024625157e-16"
$ws.Range("H2").Value = [double]"6.483054158395075e-16"
$ws.Range("K2").Value = [double]"41.35926004741992"
$ws.Range("L2").Value = "[31.41129562236177, 51.307224472478076]"
$ws.Range("M2").Value = [double]"3.574918139293004e-14"
$ws.Range("N2").Value = [double]"7.149836278586008e-14"
$ws.Range("O2").Value = [double]"1.352237078121732"
$ws.Range("P2").Value = "[1.0880791372793475, 1.6163950189641163]"
$ws.Range("S2").Value = [double]"58.56517641735933"
$ws.Range("T2").Value = "[52.59085429858379, 64.53949853613487]"
$ws.Range("W2").Value = [double]"19.72948948948988"
$ws.Range("X2").Value = [double]"18.67255255255292"
$ws.Range("Y2").Value = [double]"20.78642642642684"

# Row 3 updates
$ws.Range("E3").Value = [double]"24.23000000000035"
$ws.Range("G3").Value = [double]"3.987921104453562e-13"
$ws.Range("H3").Value = [double]"1.149632265764169e-12"
$ws.Range("K3").Value = [double]"44.64317385028245"
$ws.Range("L3").Value = "[30.299825434352456, 58.98652226621245]"
$ws.Range("M3").Value = [double]"7.424206360795438e-09"
$ws.Range("N3").Value = [double]"7.424206360795438e-09"
$ws.Range("O3").Value = [double]"0.5723422051585008"
$ws.Range("P3").Value = "[0.25786846606042335, 0.8868159442565782]"
$ws.Range("Q3").Value = [double]"0.0004428945447665367"
$ws.Range("R3").Value = [double]"0.0004428945447665367"
$ws.Range("S3").Value = [double]"60.29554044297622"
$ws.Range("T3").Value = "[52.89935373520612, 67.69172715074632]"
$ws.Range("W3").Value = [double]"22.02286286286318"
$ws.Range("X3").Value = [double]"20.81015015015045"
$ws.Range("Y3").Value = [double]"23.23557557557591"
